{"js": "// Update the date paragraph (first paragraph in the body, above the table).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateRange = paragraphs.items[0].getRange();\ndateRange.insertText(\"2023-07-27 Thursday\", Word.InsertLocation.replace);\n\n// Update the answer cells inside the table (addressed by row/column so the\n// duplicate \"57\\u00f74=14, 1\" values at different positions are each replaced\n// with the correct, distinct target value).\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst cellUpdates = [\n  { row: 0, col: 0, text: \"29\u00f77=4, 1\" }, // was \"39\u00f78=4, 7\"\n  { row: 0, col: 1, text: \"75\u00f78=9, 3\" }, // was \"57\u00f74=14, 1\"\n  { row: 0, col: 2, text: \"57\u00f77=8, 1\" }, // was \"56\u00f77=8, 0\"\n  { row: 0, col: 3, text: \"88\u00f76=14, 4\" }, // was \"28\u00f75=5, 3\"\n  { row: 0, col: 4, text: \"45\u00f78=5, 5\" }, // was \"57\u00f74=14, 1\"\n  { row: 4, col: 0, text: \"70\u00f72=35, 0\" }, // was \"91\u00f73=30, 1\"\n  { row: 4, col: 1, text: \"73\u00f76=12, 1\" }, // was \"66\u00f79=7, 3\"\n  { row: 4, col: 2, text: \"26\u00f78=3, 2\" }, // was \"35\u00f73=11, 2\"\n  { row: 4, col: 3, text: \"20\u00f77=2, 6\" }, // was \"26\u00f73=8, 2\"\n  { row: 4, col: 4, text: \"72\u00f72=36, 0\" }, // was \"40\u00f73=13, 1\"\n  { row: 8, col: 0, text: \"52\u00f75=10, 2\" }, // was \"28\u00f73=9, 1\"\n  { row: 8, col: 1, text: \"92\u00f79=10, 2\" }, // was \"42\u00f77=6, 0\"\n  { row: 8, col: 2, text: \"88\u00f74=22, 0\" }, // was \"84\u00f77=12, 0\"\n  { row: 8, col: 3, text: \"61\u00f79=6, 7\" }, // was \"80\u00f72=40, 0\"\n  { row: 8, col: 4, text: \"29\u00f77=4, 1\" }, // was \"30\u00f75=6, 0\"\n  { row: 12, col: 0, text: \"33\u00f75=6, 3\" }, // was \"19\u00f75=3, 4\"\n  { row: 12, col: 1, text: \"88\u00f74=22, 0\" }, // was \"95\u00f73=31, 2\"\n  { row: 12, col: 2, text: \"37\u00f78=4, 5\" }, // was \"30\u00f76=5, 0\"\n  { row: 12, col: 3, text: \"51\u00f77=7, 2\" }, // was \"22\u00f73=7, 1\"\n  { row: 12, col: 4, text: \"46\u00f73=15, 1\" }, // was \"80\u00f73=26, 2\"\n  { row: 16, col: 0, text: \"61\u00f74=15, 1\" }, // was \"25\u00f76=4, 1\"\n  { row: 16, col: 1, text: \"31\u00f73=10, 1\" }, // was \"58\u00f78=7, 2\"\n  { row: 16, col: 2, text: \"38\u00f76=6, 2\" }, // was \"71\u00f72=35, 1\"\n  { row: 16, col: 3, text: \"96\u00f75=19, 1\" }, // was \"77\u00f75=15, 2\"\n  { row: 16, col: 4, text: \"76\u00f77=10, 6\" }, // was \"14\u00f72=7, 0\"\n];\n\nfor (const { row, col, text } of cellUpdates) {\n  const cellRange = table.getCell(row, col).getRange();\n  cellRange.insertText(text, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date line (first paragraph in the document body, above the table).\n$d.Paragraphs(1).Range.Text = '2023-07-27 Thursday'\n\n# Update the answer cells inside the table. Cells are addressed by their\n# (row, column) position (1-based, matching Word's COM Cell(row, col))\n# instead of by old text, because a couple of the old answers repeat\n# (e.g. \"57/4=14, 1\") at different positions but must become different,\n# distinct new answers.\n$t = $d.Tables(1)\n\n$cellUpdates = @(\n    @{ Row = 1; Col = 1; Text = '29\u00f77=4, 1' }  # was '39\u00f78=4, 7'\n    @{ Row = 1; Col = 2; Text = '75\u00f78=9, 3' }  # was '57\u00f74=14, 1'\n    @{ Row = 1; Col = 3; Text = '57\u00f77=8, 1' }  # was '56\u00f77=8, 0'\n    @{ Row = 1; Col = 4; Text = '88\u00f76=14, 4' }  # was '28\u00f75=5, 3'\n    @{ Row = 1; Col = 5; Text = '45\u00f78=5, 5' }  # was '57\u00f74=14, 1'\n    @{ Row = 5; Col = 1; Text = '70\u00f72=35, 0' }  # was '91\u00f73=30, 1'\n    @{ Row = 5; Col = 2; Text = '73\u00f76=12, 1' }  # was '66\u00f79=7, 3'\n    @{ Row = 5; Col = 3; Text = '26\u00f78=3, 2' }  # was '35\u00f73=11, 2'\n    @{ Row = 5; Col = 4; Text = '20\u00f77=2, 6' }  # was '26\u00f73=8, 2'\n    @{ Row = 5; Col = 5; Text = '72\u00f72=36, 0' }  # was '40\u00f73=13, 1'\n    @{ Row = 9; Col = 1; Text = '52\u00f75=10, 2' }  # was '28\u00f73=9, 1'\n    @{ Row = 9; Col = 2; Text = '92\u00f79=10, 2' }  # was '42\u00f77=6, 0'\n    @{ Row = 9; Col = 3; Text = '88\u00f74=22, 0' }  # was '84\u00f77=12, 0'\n    @{ Row = 9; Col = 4; Text = '61\u00f79=6, 7' }  # was '80\u00f72=40, 0'\n    @{ Row = 9; Col = 5; Text = '29\u00f77=4, 1' }  # was '30\u00f75=6, 0'\n    @{ Row = 13; Col = 1; Text = '33\u00f75=6, 3' }  # was '19\u00f75=3, 4'\n    @{ Row = 13; Col = 2; Text = '88\u00f74=22, 0' }  # was '95\u00f73=31, 2'\n    @{ Row = 13; Col = 3; Text = '37\u00f78=4, 5' }  # was '30\u00f76=5, 0'\n    @{ Row = 13; Col = 4; Text = '51\u00f77=7, 2' }  # was '22\u00f73=7, 1'\n    @{ Row = 13; Col = 5; Text = '46\u00f73=15, 1' }  # was '80\u00f73=26, 2'\n    @{ Row = 17; Col = 1; Text = '61\u00f74=15, 1' }  # was '25\u00f76=4, 1'\n    @{ Row = 17; Col = 2; Text = '31\u00f73=10, 1' }  # was '58\u00f78=7, 2'\n    @{ Row = 17; Col = 3; Text = '38\u00f76=6, 2' }  # was '71\u00f72=35, 1'\n    @{ Row = 17; Col = 4; Text = '96\u00f75=19, 1' }  # was '77\u00f75=15, 2'\n    @{ Row = 17; Col = 5; Text = '76\u00f77=10, 6' }  # was '14\u00f72=7, 0'\n)\n\nforeach ($update in $cellUpdates) {\n    $t.Cell($update.Row, $update.Col).Range.Text = $update.Text\n}\n\n"}
